$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 56, pushing existing rows 56-63 down to 58-65.
$ws.Rows("56:57").Insert()

# New row 56: weekly update - Asterix, 1a (cosecha lavada), Región de Los Lagos
$ws.Cells.Item(56,1).Value = 1
$ws.Cells.Item(56,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(56,3).Value = "Arica y Parinacota"
$ws.Cells.Item(56,4).Value = 44617
$ws.Cells.Item(56,5).Value = 15
$ws.Cells.Item(56,6).Value = 100114001
$ws.Cells.Item(56,7).Value = "Papa"
$ws.Cells.Item(56,8).Value = "Asterix"
$ws.Cells.Item(56,9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(56,10).Value = 1000
$ws.Cells.Item(56,11).Value = 11000
$ws.Cells.Item(56,12).Value = 12000
$ws.Cells.Item(56,13).Value = 11500
$ws.Cells.Item(56,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(56,15).Value = "Región de Los Lagos"
$ws.Cells.Item(56,16).Value = 460
$ws.Cells.Item(56,17).Value = 25
$ws.Cells.Item(56,18).Value = "Hortaliza"

# New row 57: weekly update - Rosara, 1a (cosecha), Región del Maule
$ws.Cells.Item(57,1).Value = 1
$ws.Cells.Item(57,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57,3).Value = "Arica y Parinacota"
$ws.Cells.Item(57,4).Value = 44617
$ws.Cells.Item(57,5).Value = 15
$ws.Cells.Item(57,6).Value = 100114001
$ws.Cells.Item(57,7).Value = "Papa"
$ws.Cells.Item(57,8).Value = "Rosara"
$ws.Cells.Item(57,9).Value = "1a (cosecha)"
$ws.Cells.Item(57,10).Value = 1000
$ws.Cells.Item(57,11).Value = 9000
$ws.Cells.Item(57,12).Value = 10000
$ws.Cells.Item(57,13).Value = 9500
$ws.Cells.Item(57,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(57,15).Value = "Región del Maule"
$ws.Cells.Item(57,16).Value = 380
$ws.Cells.Item(57,17).Value = 25
$ws.Cells.Item(57,18).Value = "Hortaliza"
